$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values from 7573 to 7293 for rows 2-252
$ws.Range("C2:C252").Value = 7293
